# CIERRRE 8 OCT 2021
# Advance the payroll week from SEMANA 40 (27 Sep - 03 Oct 2021) to
# SEMANA 41 (04 Oct - 10 Oct 2021), and zero out the two pending
# "saldo"/loan cells (K4, K21) for the new week's first entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# B9 is the master "SEMANA ..." label; H9, B27, H27, B43, H43, B60 all
# reference it via formulas, so updating B9 ripples through all of them.
$ws.Range("B9").Value = "SEMANA   41  DEL    04      Al    10   DE   OCTUBRE          2021"

# Reset the two outstanding balances to 0 for the new week.
$ws.Range("K4").Value = 0
$ws.Range("K21").Value = 0

# Restore the view to where it was left (scrolled near the top, with the
# relevant cell selected) instead of the prior week's closing position.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("C21").Select()
